$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 198-199; Excel shifts the existing
# rows 198:309 down to 200:311 and carries the D-column date style along.
$ws.Rows("198:199").Insert()

# New record for row 198
$ws.Range("A198").Value = 10
$ws.Range("B198").Value = "Vega Modelo de Temuco"
$ws.Range("C198").Value = "La Araucanía"
$ws.Range("D198").Value = 44529
$ws.Range("E198").Value = 9
$ws.Range("F198").Value = 100114014
$ws.Range("G198").Value = "Betarraga"
$ws.Range("H198").Value = "Sin especificar"
$ws.Range("I198").Value = "Primera"
$ws.Range("J198").Value = 110
$ws.Range("K198").Value = 9000
$ws.Range("L198").Value = 9000
$ws.Range("M198").Value = 9000
$ws.Range("N198").Value = '$/docena de paquetes'
$ws.Range("O198").Value = "Provincia de Cautín"
$ws.Range("P198").Value = 750
$ws.Range("Q198").Value = 12
$ws.Range("R198").Value = "Hortaliza"

# New record for row 199
$ws.Range("A199").Value = 10
$ws.Range("B199").Value = "Vega Modelo de Temuco"
$ws.Range("C199").Value = "La Araucanía"
$ws.Range("D199").Value = 44529
$ws.Range("E199").Value = 9
$ws.Range("F199").Value = 100114014
$ws.Range("G199").Value = "Betarraga"
$ws.Range("H199").Value = "Sin especificar"
$ws.Range("I199").Value = "Primera"
$ws.Range("J199").Value = 250
$ws.Range("K199").Value = 700
$ws.Range("L199").Value = 700
$ws.Range("M199").Value = 700
$ws.Range("N199").Value = '$/paquete 5 unidades'
$ws.Range("O199").Value = "Región del Maule"
$ws.Range("P199").Value = 140
$ws.Range("Q199").Value = 5
$ws.Range("R199").Value = "Hortaliza"
